$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Fix D514:D516 so the bsecode values are stored as real numbers ---
# (they were previously imported as text/inline strings)
$ws.Range("D514").Value = 500124
$ws.Range("D515").Value = 524494
$ws.Range("D516").Value = 532522

# --- Append the newly scraped rows (517-525) from stock.yaml ---
$newRows = @(
    @(1, "ATUL",       "Atul Limited",                            "500027", 0.93,  7990.45, 47991,    "day", "05/09/2024 11:35:37"),
    @(2, "NESTLEIND",  "Nestle India Limited",                    "500790", -1.18, 2504.9,  1123540,  "day", "05/09/2024 11:35:37"),
    @(3, "COROMANDEL", "Coromandel International Limited",        "506395", -0.2,  1720.1,  182520,   "day", "05/09/2024 11:35:37"),
    @(4, "CHAMBLFERT", "Chambal Fertilizers & Chemicals Limited",  "500085", -3.3,  513.7,   3623113,  "day", "05/09/2024 11:35:37"),
    @(5, "EXIDEIND",   "Exide Industries Limited",                "500086", 2.38,  495.65,  4648545,  "day", "05/09/2024 11:35:37"),
    @(6, "BANKBARODA", "Bank Of Baroda",                           "532134", 0.14,  243.85,  14262866, "day", "05/09/2024 11:35:37"),
    @(7, "RBLBANK",    "Rbl Bank Limited",                         "540065", -0.43, 215.96,  4086046,  "day", "05/09/2024 11:35:37"),
    @(8, "BANDHANBNK", "Bandhan Bank Ltd",                         "541153", 2.03,  203.66,  27434822, "day", "05/09/2024 11:35:37"),
    @(9, "PNB",        "Punjab National Bank",                     "532461", 0.41,  113.4,   18547467, "day", "05/09/2024 11:35:37")
)

$startRow = 517
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    # bsecode kept as text for the freshly-added rows (not yet normalized) -
    # leading apostrophe forces Excel to store it as a text value; reset the
    # cell style afterwards so it doesn't pick up a quote-prefix format
    $ws.Cells.Item($row, 4).Value = "'" + $data[3]
    $ws.Cells.Item($row, 4).Style = "Normal"
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $ws.Cells.Item($row, 7).Value = $data[6]
    $ws.Cells.Item($row, 8).Value = $data[7]
    $ws.Cells.Item($row, 9).Value = $data[8]
}
